$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ropes"
$ws.Range("A3").Value = "Base"
$ws.Range("A4").Value = "Tough"
$ws.Range("A5").Value = "Shot"
$ws.Range("A6").Value = "Fiddle"
$ws.Range("A7").Value = "Lickety"
$ws.Range("A8").Value = " Silver"
$ws.Range("A9").Value = "Surgery"
$ws.Range("A10").Value = "Drawing"
$ws.Range("A11").Value = "Jumping"
$ws.Range("A12").Value = "Sock In It"

$ws.Range("A12").Select()
